$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# Insert a new row at position 2, shifting existing data down.
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row.
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ETH"
$ws.Range("T2").Value = 2435.5610999999999

# Update AutoFilter range to cover the new extent.
$ws.Range("A1:AB219").AutoFilter()

# Adjust the sheet view: clear frozen/top-left cell pin and move the selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("K19").Select()
